$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 46, shifting existing rows 46:147 down to 47:148.
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with its data.
$ws.Range("A46").Value = 1
$ws.Range("B46").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C46").Value = "Arica y Parinacota"
$ws.Range("D46").Value = 44987
$ws.Range("E46").Value = 15
$ws.Range("F46").Value = "Fruta"
$ws.Range("G46").Value = 100102
$ws.Range("H46").Value = "Cítricos"
$ws.Range("I46").Value = 100102004
$ws.Range("J46").Value = "Mandarina"
$ws.Range("K46").Value = "Murcott"
$ws.Range("L46").Value = "Tercera"
$ws.Range("M46").Value = 240
$ws.Range("N46").Value = 21000
$ws.Range("O46").Value = 22000
$ws.Range("P46").Value = 21417
$ws.Range("Q46").Value = "$/caja 20 kilos granel"
$ws.Range("R46").Value = "Región de O'Higgins"
$ws.Range("S46").Value = 1071
$ws.Range("T46").Value = 20
